$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 102 ("01-01-2021") with revised figures
$ws.Range("B102").Value = 38076
$ws.Range("C102").Value = 1175
$ws.Range("D102").Value = 3458
$ws.Range("E102").Value = 35813
$ws.Range("F102").Value = 575
$ws.Range("G102").Value = 581
$ws.Range("H102").Value = 35811
$ws.Range("I102").Value = 38129
$ws.Range("J102").Value = 39368

# Append new row 103 ("01-04-2021") with its figures.
# Force the date-shaped label to stay literal text (not auto-converted to
# a date serial) by switching the cell to a text format before entry, then
# restoring the default "Normal" style so the cell matches the rest of the
# column (no explicit style index).
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = "01-04-2021"
$ws.Range("A103").Style = "Normal"

$ws.Range("B103").Value = 39677
$ws.Range("C103").Value = 1242
$ws.Range("D103").Value = 3892
$ws.Range("E103").Value = 37058
$ws.Range("F103").Value = 719
$ws.Range("G103").Value = 690
$ws.Range("H103").Value = 37089
$ws.Range("I103").Value = 39882
$ws.Range("J103").Value = 41667
